$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-08-23 Saturday", $false, $false, $false, $false, $false,
                         $true, 1, $false, "2025-08-24 Sunday", 2)

# Update the division problems in the table. Each cell is addressed directly
# by (row, column) to avoid any ambiguity between old/new values that overlap
# across different cells.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "10÷8=" },
    @{ Row = 1;  Col = 2; Text = "63÷6=" },
    @{ Row = 1;  Col = 3; Text = "20÷5=" },
    @{ Row = 1;  Col = 4; Text = "96÷6=" },
    @{ Row = 1;  Col = 5; Text = "35÷4=" },

    @{ Row = 5;  Col = 1; Text = "93÷5=" },
    @{ Row = 5;  Col = 2; Text = "30÷5=" },
    @{ Row = 5;  Col = 3; Text = "80÷2=" },
    @{ Row = 5;  Col = 4; Text = "76÷8=" },
    @{ Row = 5;  Col = 5; Text = "43÷8=" },

    @{ Row = 9;  Col = 1; Text = "89÷9=" },
    @{ Row = 9;  Col = 2; Text = "95÷8=" },
    @{ Row = 9;  Col = 3; Text = "47÷2=" },
    @{ Row = 9;  Col = 4; Text = "38÷9=" },
    @{ Row = 9;  Col = 5; Text = "76÷6=" },

    @{ Row = 13; Col = 1; Text = "77÷5=" },
    @{ Row = 13; Col = 2; Text = "26÷4=" },
    @{ Row = 13; Col = 3; Text = "29÷6=" },
    @{ Row = 13; Col = 4; Text = "84÷3=" },
    @{ Row = 13; Col = 5; Text = "97÷3=" },

    @{ Row = 17; Col = 1; Text = "24÷8=" },
    @{ Row = 17; Col = 2; Text = "85÷7=" },
    @{ Row = 17; Col = 3; Text = "85÷3=" },
    @{ Row = 17; Col = 4; Text = "26÷7=" },
    @{ Row = 17; Col = 5; Text = "81÷9=" }
)

foreach ($u in $updates) {
    $t.Cell($u.Row, $u.Col).Range.Text = $u.Text
}
